$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the existing row 531. This shifts
# all rows from 531 downward (previously 531..564) down by two positions
# (they become 533..566), preserving all their values/formatting exactly.
$ws.Rows("531:532").Insert()

# Row 531 (new first "Apio - Primera" record for the new week, Coquimbo)
$ws.Range("A531").Value = 6
$ws.Range("B531").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C531").Value = "Metropolitana"
$ws.Range("D531").Value = 44610
$ws.Range("E531").Value = 13
$ws.Range("F531").Value = 100112017
$ws.Range("G531").Value = "Apio"
$ws.Range("H531").Value = "Americana (o)"
$ws.Range("I531").Value = "Primera"
$ws.Range("J531").Value = 2800
$ws.Range("K531").Value = 6000
$ws.Range("L531").Value = 7000
$ws.Range("M531").Value = 6429
$ws.Range("N531").Value = "$/docena de matas"
$ws.Range("O531").Value = "Región de Coquimbo"
$ws.Range("P531").Value = 1072
$ws.Range("Q531").Value = 6
$ws.Range("R531").Value = "Hortaliza"

# Row 532 (new first "Apio - Segunda" record for the new week, Coquimbo)
$ws.Range("A532").Value = 6
$ws.Range("B532").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C532").Value = "Metropolitana"
$ws.Range("D532").Value = 44610
$ws.Range("E532").Value = 13
$ws.Range("F532").Value = 100112017
$ws.Range("G532").Value = "Apio"
$ws.Range("H532").Value = "Americana (o)"
$ws.Range("I532").Value = "Segunda"
$ws.Range("J532").Value = 800
$ws.Range("K532").Value = 5000
$ws.Range("L532").Value = 5000
$ws.Range("M532").Value = 5000
$ws.Range("N532").Value = "$/docena de matas"
$ws.Range("O532").Value = "Región de Coquimbo"
$ws.Range("P532").Value = 833
$ws.Range("Q532").Value = 6
$ws.Range("R532").Value = "Hortaliza"
